# Applies the "All TF dc gains stable1" edit to the LMO-P2DM cell-parameter
# workbook:
#   - Renames the "Double-layer CPE frequency breakpoint" (wDL, rad/s)
#     parameter to "Double-layer CPE-integrator time constant" (tauDL, s)
#     in both the Negative- and Positive-Electrode parameter blocks, and
#     updates the associated numeric values/format.
#   - Tweaks a handful of other positive-electrode double-layer / solid-phase
#     diffusivity parameter values.
#   - Restores the active selection to match the author's final cursor spot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# --- Negative Electrode Parameters: row 36 (wDL -> tauDL) -----------------
$ws.Range("C36").Value = "tauDL"
$ws.Range("D36").Value = "\tau_\mathrm{dl}"
$ws.Range("B36").Value = "Double-layer CPE-integrator time constant"
$ws.Range("E36").Value = 1000
$ws.Range("E36").NumberFormat = "0.00E+00"
$ws.Range("G36").Value = "s"

# --- Positive Electrode Parameters: other value tweaks ---------------------
$ws.Range("E65").Value = 0.8     # nF  - Solid-phase diffusivity CPE factor
$ws.Range("E70").Value = 10      # Cdl - Double-layer capacitance
$ws.Range("E71").Value = 0.5     # nDL - Double-layer CPE factor

# --- Positive Electrode Parameters: row 72 (wDL -> tauDL) ------------------
$ws.Range("B72").Value = "Double-layer CPE-integrator time constant"
$ws.Range("C72").Value = "tauDL"
$ws.Range("D72").Value = "\tau_\mathrm{dl}"
$ws.Range("E72").Value = 100
$ws.Range("E72").NumberFormat = "0.00E+00"
$ws.Range("G72").Value = "s"

# --- Restore author's final selection (E71) ---------------------------------
$ws.Activate()
$ws.Range("E71").Select()
